$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F: width, header + data, matching header style from column E
$ws.Columns("F").ColumnWidth = 13.15

$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Country"

$ws.Range("F2").Value = "Egypt"
$ws.Range("F3").Value = "Ksa"
$ws.Range("F4").Value = "All countries"

# Data validation dropdown list on F2:F4
$ws.Range("F2:F4").Validation.Add(3, 1, 1, '"Egypt, Ksa, All countries"')

# View state: zoom + selection
$excel.ActiveWindow.Zoom = 100
$ws.Range("F6").Select()
